# Update the "Platform Coverage" sheet:
#  - the "Treatment / Campaign / MDA" row with age band 50-65 (row 5) is
#    removed, shifting the "Vaccine" rows (and the Vector Control row) up
#    by one
#  - the remaining "Treatment / Campaign / MDA" row (row 4, age band 15-50)
#    has its max-age bumped from 50 to 65

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Platform Coverage")

# Remove row 5 entirely, shifting rows 6-9 up to become rows 5-8.
$ws.Rows.Item(5).Delete() | Out-Null

# Bump the max age for the remaining MDA/Campaign row (row 4) from 50 to 65.
$ws.Range("G4").Value = 65

# Match the author's final cell selection.
$ws.Range("G12").Select() | Out-Null
